# Insert a new data row at row 844 (a new "2026/02/20" sample), pushing the
# existing rows 844-885 down to 845-886. This mirrors the diff, which shows
# every row from 844 onward shifted down by one with a brand-new row 844.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 844:885 down to 845:886, creating a blank row 844.
$ws.Rows("844:844").Insert()

# The new row's A-cell holds a literal text date (matching the sheet's
# convention of storing dates as plain strings, not Excel date serials).
# Briefly force a text number format so the assignment isn't reinterpreted
# as a date, then clear the format again so the cell carries no style -
# exactly like its sibling data cells.
$ws.Range("A844").NumberFormat = "@"
$ws.Range("A844").Value = "2026/02/20"
$ws.Range("A844").ClearFormats()

$ws.Range("B844").Value = "金"
$ws.Range("C844").Value = 4
$ws.Range("D844").Value = 201
